$wb = $excel.ActiveWorkbook

# 1. Rename header in Sheet1 (Weekly Quantity)
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# 2. Rename header in Sheet2 (Monthly Trend)
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 3. Add new sheet "PO Forecast" after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# 4. Write header row with formatting matching other sheets
$wsForecast.Cells.Item(1, 1).Value = "ds"
$wsForecast.Cells.Item(1, 2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1, 3).Value = "yhat_lower"
$wsForecast.Cells.Item(1, 4).Value = "yhat_upper"
$hdr = $wsForecast.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# 5. Write data rows
$wsForecast.Cells.Item(2, 1).Value = 45249.99999999999
$wsForecast.Cells.Item(2, 2).Value = 88
$wsForecast.Cells.Item(2, 3).Value = -385.4547155052442
$wsForecast.Cells.Item(2, 4).Value = 538.3771526228135
$wsForecast.Cells.Item(3, 1).Value = 45256.99999999999
$wsForecast.Cells.Item(3, 2).Value = 96
$wsForecast.Cells.Item(3, 3).Value = -366.8250108032742
$wsForecast.Cells.Item(3, 4).Value = 539.7298198549739
$wsForecast.Cells.Item(4, 1).Value = 45270.99999999999
$wsForecast.Cells.Item(4, 2).Value = 113
$wsForecast.Cells.Item(4, 3).Value = -324.8549650505007
$wsForecast.Cells.Item(4, 4).Value = 577.0237412591141
$wsForecast.Cells.Item(5, 1).Value = 45277.99999999999
$wsForecast.Cells.Item(5, 2).Value = 121
$wsForecast.Cells.Item(5, 3).Value = -325.366503877273
$wsForecast.Cells.Item(5, 4).Value = 599.8453559856623
$wsForecast.Cells.Item(6, 1).Value = 45298.99999999999
$wsForecast.Cells.Item(6, 2).Value = 146
$wsForecast.Cells.Item(6, 3).Value = -330.8476916209426
$wsForecast.Cells.Item(6, 4).Value = 565.6566704246741
$wsForecast.Cells.Item(7, 1).Value = 45305.99999999999
$wsForecast.Cells.Item(7, 2).Value = 154
$wsForecast.Cells.Item(7, 3).Value = -326.7026989422174
$wsForecast.Cells.Item(7, 4).Value = 617.3685295017681
$wsForecast.Cells.Item(8, 1).Value = 45312.99999999999
$wsForecast.Cells.Item(8, 2).Value = 163
$wsForecast.Cells.Item(8, 3).Value = -287.0454425700267
$wsForecast.Cells.Item(8, 4).Value = 643.0181054344242
$wsForecast.Cells.Item(9, 1).Value = 45333.99999999999
$wsForecast.Cells.Item(9, 2).Value = 187
$wsForecast.Cells.Item(9, 3).Value = -294.5379035275467
$wsForecast.Cells.Item(9, 4).Value = 659.6882839222692
$wsForecast.Cells.Item(10, 1).Value = 45340.99999999999
$wsForecast.Cells.Item(10, 2).Value = 196
$wsForecast.Cells.Item(10, 3).Value = -198.4725780143651
$wsForecast.Cells.Item(10, 4).Value = 656.1041268892556
$wsForecast.Cells.Item(11, 1).Value = 45347.99999999999
$wsForecast.Cells.Item(11, 2).Value = 204
$wsForecast.Cells.Item(11, 3).Value = -237.1972156653797
$wsForecast.Cells.Item(11, 4).Value = 652.0913919809955
$wsForecast.Cells.Item(12, 1).Value = 45354.99999999999
$wsForecast.Cells.Item(12, 2).Value = 212
$wsForecast.Cells.Item(12, 3).Value = -236.5178462160723
$wsForecast.Cells.Item(12, 4).Value = 640.9089323903196
$wsForecast.Cells.Item(13, 1).Value = 45361.99999999999
$wsForecast.Cells.Item(13, 2).Value = 221
$wsForecast.Cells.Item(13, 3).Value = -213.7392627232441
$wsForecast.Cells.Item(13, 4).Value = 691.9609590559801
$wsForecast.Cells.Item(14, 1).Value = 45368.99999999999
$wsForecast.Cells.Item(14, 2).Value = 229
$wsForecast.Cells.Item(14, 3).Value = -223.2914812039173
$wsForecast.Cells.Item(14, 4).Value = 670.398054005359
$wsForecast.Cells.Item(15, 1).Value = 45375.99999999999
$wsForecast.Cells.Item(15, 2).Value = 237
$wsForecast.Cells.Item(15, 3).Value = -202.998286652822
$wsForecast.Cells.Item(15, 4).Value = 673.4364332354985
$wsForecast.Cells.Item(16, 1).Value = 45382.99999999999
$wsForecast.Cells.Item(16, 2).Value = 245
$wsForecast.Cells.Item(16, 3).Value = -186.9303475205317
$wsForecast.Cells.Item(16, 4).Value = 667.7855128720286
$wsForecast.Cells.Item(17, 1).Value = 45389.99999999999
$wsForecast.Cells.Item(17, 2).Value = 254
$wsForecast.Cells.Item(17, 3).Value = -204.7888981387522
$wsForecast.Cells.Item(17, 4).Value = 687.9365285884893
$wsForecast.Cells.Item(18, 1).Value = 45396.99999999999
$wsForecast.Cells.Item(18, 2).Value = 262
$wsForecast.Cells.Item(18, 3).Value = -195.1250770859572
$wsForecast.Cells.Item(18, 4).Value = 723.0310652781119
$wsForecast.Cells.Item(19, 1).Value = 45403.99999999999
$wsForecast.Cells.Item(19, 2).Value = 270
$wsForecast.Cells.Item(19, 3).Value = -159.999769095825
$wsForecast.Cells.Item(19, 4).Value = 716.8140482126106
$wsForecast.Cells.Item(20, 1).Value = 45410.99999999999
$wsForecast.Cells.Item(20, 2).Value = 278
$wsForecast.Cells.Item(20, 3).Value = -165.8395766924404
$wsForecast.Cells.Item(20, 4).Value = 759.5821562696113
$wsForecast.Cells.Item(21, 1).Value = 45417.99999999999
$wsForecast.Cells.Item(21, 2).Value = 287
$wsForecast.Cells.Item(21, 3).Value = -155.8512011836277
$wsForecast.Cells.Item(21, 4).Value = 736.4542307560412
$wsForecast.Cells.Item(22, 1).Value = 45424.99999999999
$wsForecast.Cells.Item(22, 2).Value = 295
$wsForecast.Cells.Item(22, 3).Value = -129.1678202031489
$wsForecast.Cells.Item(22, 4).Value = 740.6492346510637
$wsForecast.Cells.Item(23, 1).Value = 45431.99999999999
$wsForecast.Cells.Item(23, 2).Value = 303
$wsForecast.Cells.Item(23, 3).Value = -126.477882475191
$wsForecast.Cells.Item(23, 4).Value = 763.1816964991193
$wsForecast.Cells.Item(24, 1).Value = 45438.99999999999
$wsForecast.Cells.Item(24, 2).Value = 312
$wsForecast.Cells.Item(24, 3).Value = -141.2176619128113
$wsForecast.Cells.Item(24, 4).Value = 756.3433640309426
$wsForecast.Cells.Item(25, 1).Value = 45445.99999999999
$wsForecast.Cells.Item(25, 2).Value = 320
$wsForecast.Cells.Item(25, 3).Value = -106.6519336404692
$wsForecast.Cells.Item(25, 4).Value = 740.7954719339033
$wsForecast.Cells.Item(26, 1).Value = 45459.99999999999
$wsForecast.Cells.Item(26, 2).Value = 336
$wsForecast.Cells.Item(26, 3).Value = -112.2302480747607
$wsForecast.Cells.Item(26, 4).Value = 799.8315464138835
$wsForecast.Cells.Item(27, 1).Value = 45466.99999999999
$wsForecast.Cells.Item(27, 2).Value = 345
$wsForecast.Cells.Item(27, 3).Value = -85.74848684713096
$wsForecast.Cells.Item(27, 4).Value = 790.9864389813526
$wsForecast.Cells.Item(28, 1).Value = 45501.99999999999
$wsForecast.Cells.Item(28, 2).Value = 386
$wsForecast.Cells.Item(28, 3).Value = -49.95410730917225
$wsForecast.Cells.Item(28, 4).Value = 827.2912886755003
$wsForecast.Cells.Item(29, 1).Value = 45508.99999999999
$wsForecast.Cells.Item(29, 2).Value = 394
$wsForecast.Cells.Item(29, 3).Value = -60.56578340588877
$wsForecast.Cells.Item(29, 4).Value = 865.7971421163967
$wsForecast.Cells.Item(30, 1).Value = 45522.99999999999
$wsForecast.Cells.Item(30, 2).Value = 411
$wsForecast.Cells.Item(30, 3).Value = -35.8592628626143
$wsForecast.Cells.Item(30, 4).Value = 868.4976468980345
$wsForecast.Cells.Item(31, 1).Value = 45529.99999999999
$wsForecast.Cells.Item(31, 2).Value = 419
$wsForecast.Cells.Item(31, 3).Value = -24.29949755260103
$wsForecast.Cells.Item(31, 4).Value = 872.9783881141992
$wsForecast.Cells.Item(32, 1).Value = 45536.99999999999
$wsForecast.Cells.Item(32, 2).Value = 427
$wsForecast.Cells.Item(32, 3).Value = -4.439536467235937
$wsForecast.Cells.Item(32, 4).Value = 893.8420595841549
$wsForecast.Cells.Item(33, 1).Value = 45557.99999999999
$wsForecast.Cells.Item(33, 2).Value = 452
$wsForecast.Cells.Item(33, 3).Value = 5.411920235362064
$wsForecast.Cells.Item(33, 4).Value = 881.2878728594495
$wsForecast.Cells.Item(34, 1).Value = 45564.99999999999
$wsForecast.Cells.Item(34, 2).Value = 460
$wsForecast.Cells.Item(34, 3).Value = 65.03139184363721
$wsForecast.Cells.Item(34, 4).Value = 915.976429847581
$wsForecast.Cells.Item(35, 1).Value = 45571.99999999999
$wsForecast.Cells.Item(35, 2).Value = 469
$wsForecast.Cells.Item(35, 3).Value = 49.62899230102321
$wsForecast.Cells.Item(35, 4).Value = 933.7081042650565
$wsForecast.Cells.Item(36, 1).Value = 45585.99999999999
$wsForecast.Cells.Item(36, 2).Value = 485
$wsForecast.Cells.Item(36, 3).Value = 60.73901508674183
$wsForecast.Cells.Item(36, 4).Value = 932.3056490221559
$wsForecast.Cells.Item(37, 1).Value = 45592.99999999999
$wsForecast.Cells.Item(37, 2).Value = 494
$wsForecast.Cells.Item(37, 3).Value = 3.357342208369716
$wsForecast.Cells.Item(37, 4).Value = 942.711384465164
$wsForecast.Cells.Item(38, 1).Value = 45599.99999999999
$wsForecast.Cells.Item(38, 2).Value = 502
$wsForecast.Cells.Item(38, 3).Value = 43.28532193553028
$wsForecast.Cells.Item(38, 4).Value = 951.8622755408463
$wsForecast.Cells.Item(39, 1).Value = 45634.99999999999
$wsForecast.Cells.Item(39, 2).Value = 543
$wsForecast.Cells.Item(39, 3).Value = 114.8790390479707
$wsForecast.Cells.Item(39, 4).Value = 986.2339318732783
$wsForecast.Cells.Item(40, 1).Value = 45641.99999999999
$wsForecast.Cells.Item(40, 2).Value = 551
$wsForecast.Cells.Item(40, 3).Value = 80.94047841274839
$wsForecast.Cells.Item(40, 4).Value = 1003.991901084071
$wsForecast.Cells.Item(41, 1).Value = 45648.99999999999
$wsForecast.Cells.Item(41, 2).Value = 560
$wsForecast.Cells.Item(41, 3).Value = 101.7540595730169
$wsForecast.Cells.Item(41, 4).Value = 1010.339614104617
$wsForecast.Cells.Item(42, 1).Value = 45655.99999999999
$wsForecast.Cells.Item(42, 2).Value = 568
$wsForecast.Cells.Item(42, 3).Value = 116.7185664001528
$wsForecast.Cells.Item(42, 4).Value = 979.1672762079662
$wsForecast.Cells.Item(43, 1).Value = 45662.99999999999
$wsForecast.Cells.Item(43, 2).Value = 576
$wsForecast.Cells.Item(43, 3).Value = 129.8925518705168
$wsForecast.Cells.Item(43, 4).Value = 1029.256151334001
$wsForecast.Cells.Item(44, 1).Value = 45669.99999999999
$wsForecast.Cells.Item(44, 2).Value = 585
$wsForecast.Cells.Item(44, 3).Value = 144.39213573801
$wsForecast.Cells.Item(44, 4).Value = 1046.119175033638
$wsForecast.Cells.Item(45, 1).Value = 45676.99999999999
$wsForecast.Cells.Item(45, 2).Value = 593
$wsForecast.Cells.Item(45, 3).Value = 175.2119951584875
$wsForecast.Cells.Item(45, 4).Value = 1083.051156123353
$wsForecast.Cells.Item(46, 1).Value = 45683.99999999999
$wsForecast.Cells.Item(46, 2).Value = 601
$wsForecast.Cells.Item(46, 3).Value = 144.9451847052089
$wsForecast.Cells.Item(46, 4).Value = 1017.667476516419
$wsForecast.Cells.Item(47, 1).Value = 45690.99999999999
$wsForecast.Cells.Item(47, 2).Value = 609
$wsForecast.Cells.Item(47, 3).Value = 150.3933172935496
$wsForecast.Cells.Item(47, 4).Value = 1063.82261730513

# 6. Apply date format to column A (rows 2-47)
$wsForecast.Range("A2:A47").NumberFormat = "YYYY-MM-DD HH:MM:SS"
